$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 109.666664
$ws.Range("I8").Value = 100
$ws.Range("J8").Value = 129
$ws.Range("K8").Value = 300
$ws.Range("L8").Value = 387
$ws.Range("M8").Value = -161
$ws.Range("N8").Value = -665

$ws.Range("H31").Value = 4250
$ws.Range("J31").Value = 8000
$ws.Range("L31").Value = 24000
$ws.Range("N31").Value = -24460

$ws.Range("H33").Value = 235.96
$ws.Range("I33").Value = 240.5
$ws.Range("J33").Value = 224.28572
$ws.Range("K33").Value = 240.5
$ws.Range("L33").Value = 224.28572
$ws.Range("M33").Value = -11.5
$ws.Range("N33").Value = -682.28572

$ws.Range("H38").Value = 121.90909
$ws.Range("I38").Value = 121.90909
$ws.Range("K38").Value = 365.72727
$ws.Range("M38").Value = 6.272729999999967

$ws.Range("H39").Value = 1033.2727
$ws.Range("I39").Value = 151.88889
$ws.Range("J39").Value = 4999.5
$ws.Range("K39").Value = 455.66667
$ws.Range("L39").Value = 14998.5
$ws.Range("M39").Value = -159.66667
$ws.Range("N39").Value = -15590.5

$ws.Range("H41").Value = 83803.164
$ws.Range("I41").Value = 163.4
$ws.Range("K41").Value = 163.4
$ws.Range("M41").Value = 276.6

$ws.Range("H52").Value = 143.125
$ws.Range("I52").Value = 169
$ws.Range("J52").Value = 100
$ws.Range("K52").Value = 507
$ws.Range("L52").Value = 300
$ws.Range("M52").Value = -347
$ws.Range("N52").Value = -620

$ws.Range("H54").Value = 12894.5
$ws.Range("I54").Value = 6666.6665
$ws.Range("K54").Value = 6666.6665
$ws.Range("M54").Value = -6180.6665

$ws.Range("H55").Value = 200.64706
$ws.Range("I55").Value = 207
$ws.Range("J55").Value = 185.4
$ws.Range("K55").Value = 207
$ws.Range("L55").Value = 185.4
$ws.Range("M55").Value = 7
$ws.Range("N55").Value = -613.4

$ws.Range("H64").Value = 6951.706
$ws.Range("J64").Value = 8479.727999999999
$ws.Range("L64").Value = 8479.727999999999
$ws.Range("N64").Value = -8975.727999999999

$ws.Range("H67").Value = 6951.706
$ws.Range("J67").Value = 8479.727999999999
$ws.Range("L67").Value = 8479.727999999999
$ws.Range("N67").Value = -10195.728

$ws.Range("H69").Value = 9868.691999999999
$ws.Range("I69").Value = 9124.5
$ws.Range("K69").Value = 27373.5
$ws.Range("M69").Value = -26499.5

$ws.Range("H72").Value = 9868.691999999999
$ws.Range("I72").Value = 9124.5
$ws.Range("K72").Value = 82120.5
$ws.Range("M72").Value = -77752.5

$ws.Range("H74").Value = 7375.2173
$ws.Range("I74").Value = 4672
$ws.Range("J74").Value = 8329.294
$ws.Range("K74").Value = 4672
$ws.Range("L74").Value = 8329.294
$ws.Range("M74").Value = -3736
$ws.Range("N74").Value = -10201.294

$ws.Range("H77").Value = 7375.2173
$ws.Range("I77").Value = 4672
$ws.Range("J77").Value = 8329.294
$ws.Range("K77").Value = 23360
$ws.Range("L77").Value = 41646.47
$ws.Range("M77").Value = -18680
$ws.Range("N77").Value = -51006.47

$ws.Range("H113").Value = 3983.1667
$ws.Range("I113").Value = 3599.6667
$ws.Range("K113").Value = 3599.6667
$ws.Range("M113").Value = -345.6667000000002

$ws.Range("H131").Value = 5247.1577
$ws.Range("I131").Value = 4621.143
$ws.Range("K131").Value = 13863.429
$ws.Range("M131").Value = -8823.429

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5620132.5
$ws.Range("I32").Value = 5749275.5
$ws.Range("K32").Value = 5749275.5
$ws.Range("M32").Value = -5748988.5

$ws.Range("H61").Value = 2781485.5
$ws.Range("I61").Value = 3337232.5
$ws.Range("K61").Value = 3337232.5
$ws.Range("M61").Value = -3337020.5

$ws.Range("H102").Value = 4970.357
$ws.Range("I102").Value = 4065
$ws.Range("J102").Value = 6600
$ws.Range("K102").Value = 4065
$ws.Range("L102").Value = 6600
$ws.Range("M102").Value = -2443
$ws.Range("N102").Value = -9844

$ws.Range("H132").Value = 467718.66
$ws.Range("I132").Value = 530116.75
$ws.Range("J132").Value = 10132.667
$ws.Range("K132").Value = 1590350.25
$ws.Range("L132").Value = 30398.001
$ws.Range("M132").Value = -1587820.25
$ws.Range("N132").Value = -35458.001

$ws.Range("H134").Value = 73999.5
$ws.Range("J134").Value = 73999.5
$ws.Range("L134").Value = 73999.5
$ws.Range("N134").Value = -84139.5

$ws.Range("H136").Value = 2781485.5
$ws.Range("I136").Value = 3337232.5
$ws.Range("K136").Value = 10011697.5
$ws.Range("M136").Value = -10009147.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H68").Value = 100295
$ws.Range("J68").Value = 100295
$ws.Range("L68").Value = 100295
$ws.Range("N68").Value = -101917

$ws.Range("H71").Value = 100295
$ws.Range("J71").Value = 100295
$ws.Range("L71").Value = 300885
$ws.Range("N71").Value = -308997

$ws.Range("H86").Value = 2182.0625
$ws.Range("I86").Value = 2333.111
$ws.Range("J86").Value = 1987.8572
$ws.Range("K86").Value = 2333.111
$ws.Range("L86").Value = 1987.8572
$ws.Range("M86").Value = -1210.111
$ws.Range("N86").Value = -4233.8572

$ws.Range("H89").Value = 2182.0625
$ws.Range("I89").Value = 2333.111
$ws.Range("J89").Value = 1987.8572
$ws.Range("K89").Value = 11665.555
$ws.Range("L89").Value = 9939.286
$ws.Range("M89").Value = -6049.555
$ws.Range("N89").Value = -21171.286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 700
$ws.Range("J14").Value = 700
$ws.Range("L14").Value = 700
$ws.Range("N14").Value = -1040

$ws.Range("H105").Value = 44044.25
$ws.Range("I105").Value = 58388.668
$ws.Range("K105").Value = 58388.668
$ws.Range("M105").Value = -56641.668

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 178598.14
$ws.Range("I2").Value = 208354.38
$ws.Range("K2").Value = 1250126.28
$ws.Range("M2").Value = -1250013.28

$ws.Range("H7").Value = 3727472.5
$ws.Range("J7").Value = 7000033.5
$ws.Range("L7").Value = 21000100.5
$ws.Range("N7").Value = -21000324.5

$ws.Range("H98").Value = 782.6
$ws.Range("J98").Value = 782.6
$ws.Range("L98").Value = 2347.8
$ws.Range("N98").Value = -5343.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6800.727
$ws.Range("I80").Value = 6153.9414
$ws.Range("K80").Value = 6153.9414
$ws.Range("M80").Value = -5155.9414

$ws.Range("H83").Value = 6800.727
$ws.Range("I83").Value = 6153.9414
$ws.Range("K83").Value = 30769.707
$ws.Range("M83").Value = -25777.707

$ws.Range("H132").Value = 710723.6
$ws.Range("I132").Value = 804999.2
$ws.Range("K132").Value = 2414997.6
$ws.Range("M132").Value = -2412467.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 561.4761999999999
$ws.Range("I16").Value = 595.4286
$ws.Range("J16").Value = 493.57144
$ws.Range("K16").Value = 595.4286
$ws.Range("L16").Value = 493.57144
$ws.Range("M16").Value = -425.4286
$ws.Range("N16").Value = -833.5714399999999

$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()

$ws.Range("H46").Value = 876.6818
$ws.Range("I46").Value = 849.8
$ws.Range("J46").Value = 934.2857
$ws.Range("K46").Value = 849.8
$ws.Range("L46").Value = 934.2857
$ws.Range("M46").Value = -661.8
$ws.Range("N46").Value = -1310.2857

$ws.Range("H80").Value = 60032
$ws.Range("J80").Value = 60032
$ws.Range("L80").Value = 60032
$ws.Range("N80").Value = -62278

$ws.Range("H81").Value = 49900
$ws.Range("J81").Value = 49900
$ws.Range("L81").Value = 49900
$ws.Range("N81").Value = -51896

$ws.Range("H83").Value = 60032
$ws.Range("J83").Value = 60032
$ws.Range("L83").Value = 180096
$ws.Range("N83").Value = -191328

$ws.Range("H84").Value = 49900
$ws.Range("J84").Value = 49900
$ws.Range("L84").Value = 149700
$ws.Range("N84").Value = -159684

$ws.Range("H93").Value = 2123.4285
$ws.Range("I93").Value = 1928.5
$ws.Range("J93").Value = 2269.625
$ws.Range("K93").Value = 1928.5
$ws.Range("L93").Value = 2269.625
$ws.Range("M93").Value = -680.5
$ws.Range("N93").Value = -4765.625

$ws.Range("H100").Value = 10400.2
$ws.Range("J100").Value = 18857.715
$ws.Range("L100").Value = 18857.715
$ws.Range("N100").Value = -19939.715

$ws.Range("H124").Value = 80000
$ws.Range("J124").Value = 80000
$ws.Range("L124").Value = 80000
$ws.Range("N124").Value = -89820

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2792.7
$ws.Range("I122").Value = 2404.75
$ws.Range("J122").Value = 3374.625
$ws.Range("K122").Value = 7214.25
$ws.Range("L122").Value = 10123.875
$ws.Range("M122").Value = -4764.25
$ws.Range("N122").Value = -15023.875

$ws.Range("H132").Value = 6101187.5
$ws.Range("J132").Value = 5663.3335
$ws.Range("L132").Value = 16990.0005
$ws.Range("N132").Value = -22050.0005
